$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value needs a text-prefix
# guard (apostrophe) so Excel stores it as text instead of silently coercing it
# to a number (which would strip trailing zeros / insignificant digits and change
# the stored representation, e.g. "48.00" -> 48).
$updates = @(
    @{ Cell = 'D2'; Value = '25.591.79'; TextGuard = $false }
    @{ Cell = 'E2'; Value = '  +4.72%  '; TextGuard = $false }
    @{ Cell = 'D3'; Value = '1.704.82'; TextGuard = $false }
    @{ Cell = 'E3'; Value = '  +2.99%  '; TextGuard = $false }
    @{ Cell = 'D4'; Value = '0.9995'; TextGuard = $true }
    @{ Cell = 'E4'; Value = '  -0.47%  '; TextGuard = $false }
    @{ Cell = 'D5'; Value = '329.90'; TextGuard = $true }
    @{ Cell = 'E5'; Value = '  +5.31%  '; TextGuard = $false }
    @{ Cell = 'D6'; Value = '0.9970'; TextGuard = $true }
    @{ Cell = 'E6'; Value = '  -0.36%  '; TextGuard = $false }
    @{ Cell = 'D7'; Value = '0.3672'; TextGuard = $true }
    @{ Cell = 'E7'; Value = '  +0.98%  '; TextGuard = $false }
    @{ Cell = 'D8'; Value = '48.00'; TextGuard = $true }
    @{ Cell = 'E8'; Value = '  +2.01%  '; TextGuard = $false }
    @{ Cell = 'D9'; Value = '0.3281'; TextGuard = $true }
    @{ Cell = 'E9'; Value = '  +0.18%  '; TextGuard = $false }
    @{ Cell = 'D10'; Value = '1.165'; TextGuard = $true }
    @{ Cell = 'E10'; Value = '  +3.28%  '; TextGuard = $false }
    @{ Cell = 'E11'; Value = '  +3.49%  '; TextGuard = $false }
    @{ Cell = 'D12'; Value = '0.9977'; TextGuard = $true }
    @{ Cell = 'E12'; Value = '  -0.37%  '; TextGuard = $false }
    @{ Cell = 'D13'; Value = '6.187'; TextGuard = $true }
    @{ Cell = 'E13'; Value = '  +3.97%  '; TextGuard = $false }
    @{ Cell = 'D14'; Value = '19.92'; TextGuard = $true }
    @{ Cell = 'E14'; Value = '  +2.23%  '; TextGuard = $false }
    @{ Cell = 'D15'; Value = '1.699.89'; TextGuard = $false }
    @{ Cell = 'E15'; Value = '  +2.59%  '; TextGuard = $false }
    @{ Cell = 'D16'; Value = '6.792'; TextGuard = $true }
    @{ Cell = 'E16'; Value = '  +2.66%  '; TextGuard = $false }
    @{ Cell = 'E17'; Value = '  +1.61%  '; TextGuard = $false }
    @{ Cell = 'D18'; Value = '0.06582'; TextGuard = $true }
    @{ Cell = 'E18'; Value = '  -0.34%  '; TextGuard = $false }
    @{ Cell = 'D19'; Value = '80.77'; TextGuard = $true }
    @{ Cell = 'E19'; Value = '  +3.26%  '; TextGuard = $false }
    @{ Cell = 'D20'; Value = '0.9969'; TextGuard = $true }
    @{ Cell = 'E20'; Value = '  -0.29%  '; TextGuard = $false }
    @{ Cell = 'D21'; Value = '6.036'; TextGuard = $true }
    @{ Cell = 'E21'; Value = '  +1.60%  '; TextGuard = $false }
    @{ Cell = 'D22'; Value = '16.08'; TextGuard = $true }
    @{ Cell = 'E22'; Value = '  +2.16%  '; TextGuard = $false }
    @{ Cell = 'D23'; Value = '13.04'; TextGuard = $true }
    @{ Cell = 'E23'; Value = '  +4.51%  '; TextGuard = $false }
    @{ Cell = 'D24'; Value = '25.580.61'; TextGuard = $false }
    @{ Cell = 'E24'; Value = '  +4.68%  '; TextGuard = $false }
    @{ Cell = 'D25'; Value = '2.454'; TextGuard = $true }
    @{ Cell = 'E25'; Value = '  -0.93%  '; TextGuard = $false }
    @{ Cell = 'D26'; Value = '2.463'; TextGuard = $true }
    @{ Cell = 'E26'; Value = '  +4.66%  '; TextGuard = $false }
    @{ Cell = 'D27'; Value = '149.65'; TextGuard = $true }
    @{ Cell = 'E27'; Value = '  +1.56%  '; TextGuard = $false }
    @{ Cell = 'D28'; Value = '19.10'; TextGuard = $true }
    @{ Cell = 'E28'; Value = '  +2.47%  '; TextGuard = $false }
    @{ Cell = 'D29'; Value = '1.268'; TextGuard = $true }
    @{ Cell = 'E29'; Value = '  +5.48%  '; TextGuard = $false }
    @{ Cell = 'D30'; Value = '1.890.24'; TextGuard = $false }
    @{ Cell = 'E30'; Value = '  +2.62%  '; TextGuard = $false }
    @{ Cell = 'D31'; Value = '128.13'; TextGuard = $true }
    @{ Cell = 'E31'; Value = '  +2.87%  '; TextGuard = $false }
    @{ Cell = 'D32'; Value = '4.098'; TextGuard = $true }
    @{ Cell = 'E32'; Value = '  +0.52%  '; TextGuard = $false }
    @{ Cell = 'D33'; Value = '5.950'; TextGuard = $true }
    @{ Cell = 'E33'; Value = '  +4.04%  '; TextGuard = $false }
    @{ Cell = 'D34'; Value = '0.08490'; TextGuard = $true }
    @{ Cell = 'E34'; Value = '  +0.27%  '; TextGuard = $false }
    @{ Cell = 'D35'; Value = '1.688'; TextGuard = $true }
    @{ Cell = 'E35'; Value = '  +1.30%  '; TextGuard = $false }
    @{ Cell = 'D36'; Value = '12.68'; TextGuard = $true }
    @{ Cell = 'E36'; Value = '  +2.52%  '; TextGuard = $false }
    @{ Cell = 'D37'; Value = '5.290'; TextGuard = $true }
    @{ Cell = 'E37'; Value = '  +1.03%  '; TextGuard = $false }
    @{ Cell = 'B38'; Value = 'TrustWalletToken'; TextGuard = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; TextGuard = $false }
    @{ Cell = 'D38'; Value = '1.265'; TextGuard = $true }
    @{ Cell = 'E38'; Value = '  +5.30%  '; TextGuard = $false }
    @{ Cell = 'B39'; Value = 'Hedera'; TextGuard = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; TextGuard = $false }
    @{ Cell = 'D39'; Value = '0.06192'; TextGuard = $true }
    @{ Cell = 'E39'; Value = '  +1.89%  '; TextGuard = $false }
    @{ Cell = 'B40'; Value = 'Algorand'; TextGuard = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; TextGuard = $false }
    @{ Cell = 'D40'; Value = '0.2114'; TextGuard = $true }
    @{ Cell = 'E40'; Value = '  +1.78%  '; TextGuard = $false }
    @{ Cell = 'B41'; Value = 'VeChain'; TextGuard = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; TextGuard = $false }
    @{ Cell = 'D41'; Value = '0.02256'; TextGuard = $true }
    @{ Cell = 'E41'; Value = '  +1.72%  '; TextGuard = $false }
    @{ Cell = 'D42'; Value = '8.443'; TextGuard = $true }
    @{ Cell = 'E42'; Value = '  +2.33%  '; TextGuard = $false }
    @{ Cell = 'D43'; Value = '0.6070'; TextGuard = $true }
    @{ Cell = 'E43'; Value = '  +2.17%  '; TextGuard = $false }
    @{ Cell = 'D44'; Value = '0.9968'; TextGuard = $true }
    @{ Cell = 'E44'; Value = '  -0.28%  '; TextGuard = $false }
    @{ Cell = 'D45'; Value = '14.02'; TextGuard = $true }
    @{ Cell = 'E45'; Value = '  +9.83%  '; TextGuard = $false }
    @{ Cell = 'D46'; Value = '3.826'; TextGuard = $true }
    @{ Cell = 'E46'; Value = '  +0.85%  '; TextGuard = $false }
    @{ Cell = 'D47'; Value = '0.5821'; TextGuard = $true }
    @{ Cell = 'E47'; Value = '  +3.01%  '; TextGuard = $false }
    @{ Cell = 'D48'; Value = '125.47'; TextGuard = $true }
    @{ Cell = 'E48'; Value = '  +2.24%  '; TextGuard = $false }
    @{ Cell = 'D49'; Value = '1.995'; TextGuard = $true }
    @{ Cell = 'E49'; Value = '  +2.23%  '; TextGuard = $false }
    @{ Cell = 'D50'; Value = '0.07198'; TextGuard = $true }
    @{ Cell = 'E50'; Value = '  +4.16%  '; TextGuard = $false }
    @{ Cell = 'D51'; Value = '1.204'; TextGuard = $true }
    @{ Cell = 'E51'; Value = '  +3.24%  '; TextGuard = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.TextGuard) {
        $cell.Value = "'" + $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
